# Natmi following Dr Hou advice
# Expand the Has2->Hmmr ligand-receptor pair table (Sheet1) from 2 data rows
# to the full 6-row sending/target cluster combination table (ECs, FAPs, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Has2", "Hmmr", "ECs", 2, 0.6666666666666666, 0.3897156666666666, 1.169147, 0.01604749516233658, 0.01604749516233658, 3, 1, 2.161741333333333, 6.485224, 0.903256933110598, 0.903256933110598, 0.8424644648808887, 7.582180183927999, 0.0144950112644393, 0.0144950112644393),
    @("ECs", "Has2", "Hmmr", "FAPs", 2, 0.6666666666666666, 0.3897156666666666, 1.169147, 0.01604749516233658, 0.01604749516233658, 2, 0.6666666666666666, 0.2315326666666667, 0.694598, 0.09674306688940201, 0.096743066889402, 0.09023190754511111, 0.812087167906, 0.001552483897897283, 0.001552483897897283),
    @("FAPs", "Has2", "Hmmr", "ECs", 3, 1, 23.535792, 70.607376, 0.9691437644584301, 0.9691437644584302, 3, 1, 2.161741333333333, 6.485224, 0.903256933110598, 0.903256933110598, 50.87829437913599, 457.904649412224, 0.8753858244279813, 0.8753858244279814),
    @("FAPs", "Has2", "Hmmr", "FAPs", 3, 1, 23.535792, 70.607376, 0.9691437644584301, 0.9691437644584302, 2, 0.6666666666666666, 0.2315326666666667, 0.694598, 0.09674306688940201, 0.096743066889402, 5.449304683872001, 49.043742154848, 0.09375794003044877, 0.09375794003044877),
    @("sCs", "Has2", "Hmmr", "ECs", 3, 1, 0.3596323333333333, 1.078897, 0.01480874037923328, 0.01480874037923328, 3, 1, 2.161741333333333, 6.485224, 0.903256933110598, 0.903256933110598, 0.7774320797697777, 6.996888717928, 0.01337609741817733, 0.01337609741817733),
    @("sCs", "Has2", "Hmmr", "FAPs", 3, 1, 0.3596323333333333, 1.078897, 0.01480874037923328, 0.01480874037923328, 2, 0.6666666666666666, 0.2315326666666667, 0.694598, 0.09674306688940201, 0.096743066889402, 0.08326663315622224, 0.7493996984060001, 0.001432642961055954, 0.001432642961055954)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = 2 + $i
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
